# Apply the "update with new logo and colors" metadata refresh to the
# RxNorm - Fludarabine ValueSet workbook:
#   - bump Version 0.1.6 -> 0.1.7
#   - Status active -> draft
#   - Date refreshed
#   - Publisher Contact now shows the org display text plus a second
#     "Contact" row for the named maintainer (Bob Milius)
#   - a new "Jurisdiction" property row is added
# The "Include from RxNorm" sheet content is unaffected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# --- Update existing scalar properties in place -----------------------
$ws1.Range("B3").Value  = "0.1.7"
$ws1.Range("B6").Value  = "draft"
$ws1.Range("B8").Value  = "2024-08-23T10:17:11-05:00"
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- Make room for the second Contact row + the new Jurisdiction row --
# Before:  row11 = Contact/"No display for ContactDetail", row12 = Description, ...
# After:   row11 = Contact/Bob Milius, row12 = Jurisdiction, row13 = Description, ...
$ws1.Rows.Item(12).Insert()

$ws1.Range("A11").Value = "Contact"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""
